$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 14
$ws.Range("B3").Value = "testdriver@gmail.com"
$ws.Range("C3").Value = "Test "
$ws.Range("D3").Value = "Driver"
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Driver"
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = $false
